# Leave Card update 12/22/2023 10:59 AM
# Updates the "2018 LEAVE CREDITS" sheet: fixes the PERIOD dates for 2023
# (rows 18-29) to explicit literal dates (no more EDATE()/shared formulas),
# and fills in the EARNED (C) column for Apr-Sep 2023 (rows 21-26) with 1.25
# each. Dependent formulas (G21:G26 "EARNED ", E9/I9 BALANCE totals, and the
# CONVERTION sheet A7 total) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# --- PERIOD column (A) : replace formulas with literal serial dates ---
$ws.Cells.Item(18, 1).Value = 44957
$ws.Cells.Item(19, 1).Value = 44985
$ws.Cells.Item(20, 1).Value = 45016
$ws.Cells.Item(21, 1).Value = 45046
$ws.Cells.Item(22, 1).Value = 45077
$ws.Cells.Item(23, 1).Value = 45107
$ws.Cells.Item(24, 1).Value = 45138
$ws.Cells.Item(25, 1).Value = 45169
$ws.Cells.Item(26, 1).Value = 45199
$ws.Cells.Item(27, 1).Value = 45230
$ws.Cells.Item(28, 1).Value = 45260
$ws.Cells.Item(29, 1).Value = 45291

# --- EARNED column (C) : fill in 1.25 for Apr-Sep 2023 (rows 21-26) ---
$ws.Cells.Item(21, 3).Value = 1.25
$ws.Cells.Item(22, 3).Value = 1.25
$ws.Cells.Item(23, 3).Value = 1.25
$ws.Cells.Item(24, 3).Value = 1.25
$ws.Cells.Item(25, 3).Value = 1.25
$ws.Cells.Item(26, 3).Value = 1.25

# --- Make "2018 LEAVE CREDITS" the active/selected sheet & selection ---
$ws.Activate()
$ws.Range("C20:C26").Select()
$excel.ActiveWindow.ActivePane.TopLeftCell = $ws.Range("A21")
